# Update the "取得日時" (acquisition timestamp) column on the "ランサーズ" sheet.
# All data rows (2-14) currently show 2026-01-08 01:27:20 and must be
# refreshed to the new run's timestamp: 2026-01-08 01:59:59.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-08 01:59:59"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value -ne $null -and $cell.Value -ne "") {
        $cell.Value = $newTimestamp
    }
}
